$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported; insert it as a new row 26 and
# push the existing rows 26-69 down to 27-70 (matching the target diff).
$ws.Rows(26).Insert()

$ws.Cells.Item(26, 1).Value = 10
$ws.Cells.Item(26, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(26, 3).Value = "La Araucanía"
$ws.Cells.Item(26, 4).Value = 45002
$ws.Cells.Item(26, 5).Value = 9
$ws.Cells.Item(26, 6).Value = 100112042
$ws.Cells.Item(26, 7).Value = "Locoto"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 90
$ws.Cells.Item(26, 11).Value = 4400
$ws.Cells.Item(26, 12).Value = 4400
$ws.Cells.Item(26, 13).Value = 4400
$ws.Cells.Item(26, 14).Value = "$/kilo"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 4400
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
